# Update "want-to-go count" (F column) values across all sheets
# per the source data refresh (gh-pages output at commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1014
$ws.Range("F3").Value = 578
$ws.Range("F4").Value = 9069
$ws.Range("F8").Value = 6295
$ws.Range("F12").Value = 9332
$ws.Range("F13").Value = 10838
$ws.Range("F14").Value = 1217
$ws.Range("F15").Value = 1109
$ws.Range("F16").Value = 4864
$ws.Range("F17").Value = 778
$ws.Range("F18").Value = 425
$ws.Range("F22").Value = 1318
$ws.Range("F23").Value = 223
$ws.Range("F24").Value = 1846
$ws.Range("F25").Value = 862
$ws.Range("F26").Value = 1195
$ws.Range("F28").Value = 2005
$ws.Range("F29").Value = 410
$ws.Range("F30").Value = 595
$ws.Range("F31").Value = 2605
$ws.Range("F33").Value = 177
$ws.Range("F34").Value = 1696
$ws.Range("F36").Value = 1327
$ws.Range("F37").Value = 428
$ws.Range("F38").Value = 14
$ws.Range("F39").Value = 901
$ws.Range("F41").Value = 3267
$ws.Range("F43").Value = 79
$ws.Range("F44").Value = 494
$ws.Range("F45").Value = 567
$ws.Range("F47").Value = 891
$ws.Range("F49").Value = 4190
# Sheet: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 10
$ws.Range("F17").Value = 5
$ws.Range("F18").Value = 10
$ws.Range("F19").Value = 8
$ws.Range("F20").Value = 5
$ws.Range("F24").Value = 103
# Sheet: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5786
# Sheet: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1014
$ws.Range("F3").Value = 578
$ws.Range("F4").Value = 9069
$ws.Range("F6").Value = 10
$ws.Range("F8").Value = 6295
$ws.Range("F10").Value = 9332
$ws.Range("F11").Value = 9332
$ws.Range("F12").Value = 10838
$ws.Range("F14").Value = 1217
$ws.Range("F15").Value = 1109
$ws.Range("F16").Value = 4864
$ws.Range("F17").Value = 778
$ws.Range("F18").Value = 425
$ws.Range("F22").Value = 1318
$ws.Range("F23").Value = 223
$ws.Range("F24").Value = 862
$ws.Range("F25").Value = 1195
$ws.Range("F28").Value = 2005
$ws.Range("F29").Value = 410
$ws.Range("F30").Value = 2605
$ws.Range("F31").Value = 177
$ws.Range("F32").Value = 1696
$ws.Range("F35").Value = 428
$ws.Range("F36").Value = 10
$ws.Range("F37").Value = 8
$ws.Range("F39").Value = 901
$ws.Range("F43").Value = 79
$ws.Range("F44").Value = 494
$ws.Range("F45").Value = 567
$ws.Range("F46").Value = 891
$ws.Range("F48").Value = 4190
